# model_testing.xlsx edit — add a second "num_epochs / batch_size / MAPE / RMSE"
# summary block (columns Q:T) mirroring the existing H/I + D/E (epoch-1) data,
# and tidy up a handful of stray duplicate "applied-but-empty-fill" styles that
# Excel collapses whenever a cell's format is re-pasted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New headers Q2:T2  (num_epochs, batch_size, MAPE, RMSE)
#    Pull the formatting from the matching existing headers so the new cells
#    land on the same styles (H2/I2 = header style "3", N2/O2 = header style "1").
# ---------------------------------------------------------------------------
$ws.Range("H2:I2").Copy()
$ws.Range("Q2:R2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("N2:O2").Copy()
$ws.Range("S2:T2").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("Q2").Value = "num_epochs"
$ws.Range("R2").Value = "batch_size"
$ws.Range("S2").Value = "MAPE"
$ws.Range("T2").Value = "RMSE"

# ---------------------------------------------------------------------------
# 2) New row 3 (Q3:T3) — values + formats copied from H3:I3 (num_epochs,
#    batch_size) and D3:E3 (MAPE, RMSE), i.e. run #1's epoch/batch config and
#    its MAPE/RMSE result.
# ---------------------------------------------------------------------------
$ws.Range("H3:I3").Copy()
$ws.Range("Q3:R3").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("H3:I3").Copy()
$ws.Range("Q3:R3").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("D3:E3").Copy()
$ws.Range("S3:T3").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("D3:E3").Copy()
$ws.Range("S3:T3").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# 3) Rows 4:12 — blank placeholder cells under the new Q:T block, carrying the
#    same per-row styling as the H/I and N/O columns immediately to their left.
# ---------------------------------------------------------------------------
for ($r = 4; $r -le 12; $r++) {
    $ws.Range("H$r`:I$r").Copy()
    $ws.Range("Q$r`:R$r").PasteSpecial(-4122)   # xlPasteFormats
    $ws.Range("N$r`:O$r").Copy()
    $ws.Range("S$r`:T$r").PasteSpecial(-4122)   # xlPasteFormats
}

# ---------------------------------------------------------------------------
# 4) New column widths for Q:T (autofit to the new header/content widths).
# ---------------------------------------------------------------------------
$ws.Columns.Item(17).AutoFit()   # Q
$ws.Columns.Item(18).AutoFit()   # R
$ws.Columns.Item(20).AutoFit()   # T

# ---------------------------------------------------------------------------
# 5) Housekeeping: re-apply formats on a few cells that carried a redundant
#    "apply fill / no colour" style so they settle back onto the plain style
#    (cosmetically identical, just a tidier style table).
# ---------------------------------------------------------------------------
$ws.Range("H16").Copy()
$ws.Range("M16").PasteSpecial(-4122)
$ws.Range("M17").PasteSpecial(-4122)
$ws.Range("M19").PasteSpecial(-4122)
$ws.Range("M20").PasteSpecial(-4122)
$ws.Range("M22").PasteSpecial(-4122)
$ws.Range("M27").PasteSpecial(-4122)
$ws.Range("M28").PasteSpecial(-4122)

$ws.Range("O16").Copy()
$ws.Range("O19").PasteSpecial(-4122)

$ws.Range("H20:L20").Copy()
$ws.Range("H21:L21").PasteSpecial(-4122)

$ws.Range("M20").Copy()
$ws.Range("M21").PasteSpecial(-4122)

$ws.Range("N20").Copy()
$ws.Range("N21").PasteSpecial(-4122)

$ws.Range("O17").Copy()
$ws.Range("O21").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 6) View state: scroll so column D is the leftmost visible column, and move
#    the selection to the new block.
# ---------------------------------------------------------------------------
$ws.Range("D1").Select()
$ws.Range("S12").Select()
